{"js": "// The document contains several paragraphs where an \"<id>...</id>\" marker\n// was typed across multiple separately-formatted runs (e.g. \"<id>\", the\n// page id text, and \"</id>\" as distinct runs). The edit re-types each one\n// as a single contiguous run (same visible text), which Word collapses\n// into one run using the formatting of the first run in the matched range.\nconst body = context.document.body;\n\nconst ids = [\"p060v_1\", \"p060v_2\", \"p060v_3\", \"p060v_4\", \"p060v_5\"];\n\nfor (const id of ids) {\n  const searchText = \"<id>\" + id + \"</id>\";\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  const range = results.items[0];\n  // Re-insert the same text, replacing the matched range. This merges the\n  // previously-split runs into a single run (picking up the formatting of\n  // the first original run), matching the target OOXML.\n  range.insertText(searchText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document contains several paragraphs where an \"<id>...</id>\" marker\n# was typed across multiple separately-formatted runs (e.g. \"<id>\", the\n# page id text, and \"</id>\" as distinct runs - sometimes split even\n# further). The edit re-types each one as a single contiguous run (same\n# visible text), which Word collapses into one run using the formatting\n# of the first run in the matched range.\n$d = $word.ActiveDocument\n\n$ids = @(\"p060v_1\", \"p060v_2\", \"p060v_3\", \"p060v_4\", \"p060v_5\")\n\nforeach ($id in $ids) {\n    $text = \"<id>$id</id>\"\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # Find-and-replace the marker text with itself. Word re-creates the\n    # matched range as a single run (using the first original run's\n    # formatting), merging the previously split runs - matching the\n    # target OOXML.\n    $find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null\n}\n"}
